$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Cells.Item(3, 3).Value = 5.126259292976878
$ws.Cells.Item(3, 4).Value = 5.126259292976878
$ws.Cells.Item(3, 5).Value = 6.388817460606935
$ws.Cells.Item(3, 6).Value = 5.001228578514028
$ws.Cells.Item(3, 7).Value = 6.063222416823513
$ws.Cells.Item(3, 8).Value = 4.775233946660782
$ws.Cells.Item(3, 9).Value = 5.853006762425328
$ws.Cells.Item(3, 10).Value = 4.614202701995725
$ws.Cells.Item(3, 11).Value = 5.677967295125954
$ws.Cells.Item(3, 12).Value = 4.485023965340245
$ws.Cells.Item(3, 13).Value = 5.589650812291744
# Row 8
$ws.Cells.Item(8, 3).Value = 66.02392186623553
$ws.Cells.Item(8, 4).Value = 58.302933315735636
$ws.Cells.Item(8, 5).Value = 78.64476853532187
$ws.Cells.Item(8, 6).Value = 106.71243140319352
$ws.Cells.Item(8, 7).Value = 11.648441786877179
$ws.Cells.Item(8, 8).Value = 35.223363588025556
$ws.Cells.Item(8, 9).Value = 28.928214920170046
$ws.Cells.Item(8, 10).Value = 55.15209498655532
$ws.Cells.Item(8, 11).Value = 67.65105990529806
$ws.Cells.Item(8, 12).Value = 79.58953506712467
$ws.Cells.Item(8, 13).Value = 53.6080626987703
# Row 10
$ws.Cells.Item(10, 3).Value = 84.31863285861598
$ws.Cells.Item(10, 4).Value = 86.44920509461834
$ws.Cells.Item(10, 5).Value = 92.34812785596833
$ws.Cells.Item(10, 6).Value = 80.4180765552256
$ws.Cells.Item(10, 7).Value = 75.70309992329442
$ws.Cells.Item(10, 8).Value = 74.30334405716413
$ws.Cells.Item(10, 9).Value = 72.49620834553139
$ws.Cells.Item(10, 10).Value = 71.6548532382814
$ws.Cells.Item(10, 11).Value = 77.17073131120219
$ws.Cells.Item(10, 12).Value = 70.80233129353206
$ws.Cells.Item(10, 13).Value = 71.9411847150482
# Row 11
$ws.Cells.Item(11, 3).Value = 93.28266885373951
$ws.Cells.Item(11, 4).Value = 95.63974536223392
$ws.Cells.Item(11, 5).Value = 102.16579129076352
$ws.Cells.Item(11, 6).Value = 88.96743871480895
$ws.Cells.Item(11, 7).Value = 92.33381200292212
$ws.Cells.Item(11, 8).Value = 90.54841736231234
$ws.Cells.Item(11, 9).Value = 88.09909516039598
$ws.Cells.Item(11, 10).Value = 87.24141904605743
$ws.Cells.Item(11, 11).Value = 97.61010982043769
$ws.Cells.Item(11, 12).Value = 86.8465327534069
$ws.Cells.Item(11, 13).Value = 88.93994145595246
# Row 12
$ws.Cells.Item(12, 3).Value = 383.98971969746793
$ws.Cells.Item(12, 4).Value = 393.6924108297709
$ws.Cells.Item(12, 5).Value = 420.556285520114
$ws.Cells.Item(12, 6).Value = 366.2264648716859
$ws.Cells.Item(12, 7).Value = 360.3408664462942
$ws.Cells.Item(12, 8).Value = 353.5362380975382
$ws.Cells.Item(12, 9).Value = 344.4891350327873
$ws.Cells.Item(12, 10).Value = 340.79037573065307
$ws.Cells.Item(12, 11).Value = 373.65792364301353
$ws.Cells.Item(12, 12).Value = 337.9036525136813
$ws.Cells.Item(12, 13).Value = 344.60368168305354
# Row 17
$ws.Cells.Item(17, 3).Value = 148.73258204571576
$ws.Cells.Item(17, 4).Value = 168.4105158983527
$ws.Cells.Item(17, 5).Value = 172.191377218746
$ws.Cells.Item(17, 6).Value = 182.29494722351131
$ws.Cells.Item(17, 7).Value = 176.485134923004
$ws.Cells.Item(17, 8).Value = 171.53281008509637
$ws.Cells.Item(17, 9).Value = 162.9652308192325
$ws.Cells.Item(17, 10).Value = 163.35175286197034
$ws.Cells.Item(17, 11).Value = 154.7285151187652
$ws.Cells.Item(17, 12).Value = 165.37516340826727
$ws.Cells.Item(17, 13).Value = 163.98558897770232
# Row 18
$ws.Cells.Item(18, 3).Value = 236.40176579265258
$ws.Cells.Item(18, 4).Value = 188.7185035173737
$ws.Cells.Item(18, 5).Value = 223.6084381439636
$ws.Cells.Item(18, 6).Value = 161.432634364715
$ws.Cells.Item(18, 7).Value = 199.19716806847006
$ws.Cells.Item(18, 8).Value = 194.34088979899838
$ws.Cells.Item(18, 9).Value = 183.64923157747762
$ws.Cells.Item(18, 10).Value = 185.53353004511973
$ws.Cells.Item(18, 11).Value = 168.83289307830813
$ws.Cells.Item(18, 12).Value = 198.62804003912734
$ws.Cells.Item(18, 13).Value = 198.54212487166123
# Row 19
$ws.Cells.Item(19, 3).Value = 8.119159100686034
$ws.Cells.Item(19, 4).Value = 8.410940294790228
$ws.Cells.Item(19, 5).Value = 8.119159100686034
$ws.Cells.Item(19, 6).Value = 8.470599669657348
$ws.Cells.Item(19, 7).Value = 7.705380184764197
$ws.Cells.Item(19, 8).Value = 8.08783171097124
$ws.Cells.Item(19, 9).Value = 7.438229909453067
$ws.Cells.Item(19, 10).Value = 7.815092485708791
$ws.Cells.Item(19, 11).Value = 7.215782908475841
$ws.Cells.Item(19, 12).Value = 7.596301973165214
$ws.Cells.Item(19, 13).Value = 7.10354686796204
